$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid numeric auto-conversion of values
# like "1.00", "32.00", "0.0000172" while keeping the default cell style afterwards).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D10", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D26", "D28", "D29", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.166.81'
$ws.Range("E2").Value = '  +2.53%  '

$ws.Range("D3").Value = '3.469.31'
$ws.Range("E3").Value = '  +2.25%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '579.12'
$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("D6").Value = '148.17'
$ws.Range("E6").Value = '  +3.66%  '

$ws.Range("D7").Value = '3.469.51'
$ws.Range("E7").Value = '  +2.30%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +1.38%  '

$ws.Range("D10").Value = '7.66'
$ws.Range("E10").Value = '  +0.57%  '

$ws.Range("E11").Value = '  +2.05%  '

$ws.Range("D12").Value = '0.405'
$ws.Range("E12").Value = '  +4.98%  '

$ws.Range("D13").Value = '4.063.80'
$ws.Range("E13").Value = '  +2.21%  '

$ws.Range("D14").Value = '29.78'
$ws.Range("E14").Value = '  +6.30%  '

$ws.Range("E15").Value = '  +2.74%  '

$ws.Range("D16").Value = '3.474.73'
$ws.Range("E16").Value = '  +2.26%  '

$ws.Range("D17").Value = '0.0000172'
$ws.Range("E17").Value = '  +1.08%  '

$ws.Range("D18").Value = '63.179.19'
$ws.Range("E18").Value = '  +2.48%  '

$ws.Range("D19").Value = '6.34'
$ws.Range("E19").Value = '  +3.38%  '

$ws.Range("D20").Value = '14.43'
$ws.Range("E20").Value = '  +5.75%  '

$ws.Range("E21").Value = '  +1.12%  '

$ws.Range("D22").Value = '388.79'
$ws.Range("E22").Value = '  +0.45%  '

$ws.Range("D23").Value = '0.559'
$ws.Range("E23").Value = '  +1.96%  '

$ws.Range("D24").Value = '74.74'
$ws.Range("E24").Value = '  +0.40%  '

$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("D26").Value = '3.610.73'
$ws.Range("E26").Value = '  +2.20%  '

$ws.Range("E27").Value = '  +1.39%  '

$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  -0.90%  '

$ws.Range("D29").Value = '7.61'
$ws.Range("E29").Value = '  +3.12%  '

$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("E31").Value = '  +2.34%  '

$ws.Range("E32").Value = '  -0.68%  '

$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("E34").Value = '  -1.12%  '

$ws.Range("D35").Value = '23.67'
$ws.Range("E35").Value = '  +1.50%  '

$ws.Range("D36").Value = '5.32'
$ws.Range("E36").Value = '  +4.32%  '

$ws.Range("D37").Value = '7.09'
$ws.Range("E37").Value = '  +2.42%  '

$ws.Range("D38").Value = '32.00'
$ws.Range("E38").Value = '  +16.65%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").Value = '  +6.24%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '169.75'
$ws.Range("E40").Value = '  +0.49%  '

$ws.Range("D41").Value = '3.506.36'
$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("D42").Value = '0.0759'
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("E43").Value = '  +2.09%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = '1.22'
$ws.Range("E44").Value = '  +5.32%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '42.26'
$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.73'
$ws.Range("E46").Value = '  +3.64%  '

$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '4.44'
$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("D48").Value = '2.620.70'
$ws.Range("E48").Value = '  +5.42%  '

$ws.Range("D49").Value = '2.27'
$ws.Range("E49").Value = '  +11.74%  '

$ws.Range("D50").Value = '23.14'
$ws.Range("E50").Value = '  +1.64%  '

$ws.Range("D51").Value = '6.77'
$ws.Range("E51").Value = '  +2.05%  '

# Drop the temporary text number-format again so the cells end up with no explicit
# style (matching their original unstyled state) while keeping the stored text type.
foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
